$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Ie7b5PQ8"
$ws.Range("C3").Value = "07:00"
$ws.Range("D3").Value = "SINGAPORE - PREMIER LEAGUE"
$ws.Range("E3").Value = "Hougang"
$ws.Range("F3").Value = "Geylang"

$ws.Range("G3").Value = 4.1
$ws.Range("H3").Value = 4.5
$ws.Range("I3").Value = 1.53
$ws.Range("J3").Value = 4.33
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 34
$ws.Range("O3").Value = 1.05
$ws.Range("P3").Value = 11
$ws.Range("Q3").Value = 1.2
$ws.Range("R3").Value = 4.33
$ws.Range("S3").Value = 1.14
$ws.Range("T3").Value = 5.5
$ws.Range("U3").Value = 1.29
$ws.Range("V3").Value = 3.5
$ws.Range("W3").Value = 34
$ws.Range("X3").Value = 41
$ws.Range("Y3").Value = 17
$ws.Range("Z3").Value = 51
$ws.Range("AA3").Value = 29
$ws.Range("AB3").Value = 23
$ws.Range("AC3").Value = 34
$ws.Range("AD3").Value = 13
$ws.Range("AE3").Value = 12
$ws.Range("AF3").Value = 23
$ws.Range("AG3").Value = 51
$ws.Range("AH3").Value = 21
$ws.Range("AI3").Value = 17
$ws.Range("AJ3").Value = 10
$ws.Range("AK3").Value = 17
$ws.Range("AL3").Value = 11
$ws.Range("AM3").Value = 15
$ws.Range("AN3").Value = 8
$ws.Range("AO3").Value = 21
$ws.Range("AP3").Value = 19
$ws.Range("AQ3").Value = 51
$ws.Range("AR3").Value = 51
$ws.Range("AS3").Value = 67
$ws.Range("AT3").Value = 5.5
$ws.Range("AU3").Value = 7
$ws.Range("AV3").Value = 29
$ws.Range("AW3").Value = 4.75
$ws.Range("AX3").Value = 8
$ws.Range("AY3").Value = 11
$ws.Range("AZ3").Value = 19
$ws.Range("BA3").Value = 26
$ws.Range("BB3").Value = 51
$ws.Range("BC3").Value = 126
